# Auto-generated edit script: updates cryptos list Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.315.94'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '1.566.79'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''210.73'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''44.41'
$ws.Range("E8").Value = '  -4.01%  '
$ws.Range("D9").Value = '''23.61'
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("D10").Value = '''0.244'
$ws.Range("E10").Value = '  -1.27%  '
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("D12").Value = '''0.0894'
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("D13").Value = '1.791.18'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '1.574.70'
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '28.319.84'
$ws.Range("E16").Value = '  -0.84%  '
$ws.Range("D17").Value = '''0.513'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '''61.03'
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("D19").Value = '''227.48'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '''7.38'
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").Value = '0.0₃0676'
$ws.Range("E21").Value = '  -2.39%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '''3.92'
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("D24").Value = '''8.93'
$ws.Range("E24").Value = '  -2.20%  '
$ws.Range("E25").Value = '  -1.04%  '
$ws.Range("D26").Value = '''150.52'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("D28").Value = '''0.103'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("E29").Value = '  -1.54%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("E32").Value = '  -2.32%  '
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("D35").Value = '1.379.14'
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("D37").Value = '''1.49'
$ws.Range("E37").Value = '  -2.80%  '
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("E39").Value = '  +2.21%  '
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("D42").Value = '''1.92'
$ws.Range("E42").Value = '  +3.45%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '''0.0475'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").Value = '''0.782'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("D47").Value = '''62.12'
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("D48").Value = '''0.917'
$ws.Range("E48").Value = '  -6.19%  '
$ws.Range("D49").Value = '1.702.98'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").Value = '''85.44'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("D51").Value = '0.0₆0101'
$ws.Range("E51").Value = '  -2.00%  '
